$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 3 mirroring row 2's layout: a date/time value in A (same
# style as A2), mostly-zero sentiment metrics in B:M, and the "Random"
# method label in N (reusing the existing shared string).

# Copy A2's format (date/time number format) onto A3 before setting its value.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A3").Value = 42605.648344907408

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Random"
